# UC 12 - Invito Utenti Playlist Collaborativa
# Applies the authored edit described by the commit "upload:use cases 9 to 20S":
#   1. Remove the stray _GoBack bookmark that sat after "18" in the date line.
#   2. Rewrite the "L’interfaccia permette..." bullet into "L’utente digita...".
#   3. Split the "L’utente clicca..." bullet into two bullets: one describing
#      what the system shows, and a new one (with a fresh _GoBack bookmark)
#      describing the user clicking the desired profile.
#   4. Reword the exit-condition bullet.

$d = $word.ActiveDocument

# 1) Drop the old _GoBack bookmark near the date ("18/10/2022").
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) "L’interfaccia permette di specificare l’alias dell’utente da invitare."
#    -> "L’utente digita l’alias dell’utente da invitare nel menu di ricerca."
$null = $d.Content.Find.Execute(
  "L’interfaccia permette di specificare l’alias dell’utente da invitare.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "L’utente digita l’alias dell’utente da invitare nel menu di ricerca.", 2)

# 3) Split bullet:
#    "L’utente clicca sul profilo relativo all’utente desiderato."
#    becomes the "Il sistema mostra..." bullet, followed by a brand-new
#    bullet "L’utente clicca sul profilo relativo all’utente che intende
#    invitare." (with a _GoBack bookmark inserted between "ch" and "e").
$null = $d.Content.Find.Execute(
  "L’utente clicca sul profilo relativo all’utente desiderato.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Il sistema mostra l’elenco di utenti con il medesimo alias cercato, mostra inoltre le rispettive foto profilo e email per individuare facilmente l’utente desiderato.",
  2)

$rNewPara = $d.Content
$null = $rNewPara.Find.Execute("Il sistema mostra l’elenco di utenti con il medesimo alias cercato, mostra inoltre le rispettive foto profilo e email per individuare facilmente l’utente desiderato.")
$rNewPara.Collapse(0)
$rNewPara.InsertParagraphAfter()
$rNewPara.Collapse(0)
$rNewPara.Move(1, 1) | Out-Null
$rNewPara.InsertAfter("L’utente clicca sul profilo relativo all’utente che intende invitare.")

# Re-locate the "...utente ch" boundary now that the run has been rebuilt, and
# drop a fresh _GoBack bookmark exactly between "ch" and "e intende invitare".
$rBookmarkPoint = $d.Content
$null = $rBookmarkPoint.Find.Execute("L’utente clicca sul profilo relativo all’utente ch")
$rBookmarkPoint.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rBookmarkPoint)

# 4) "L’utente selezionato ha ora i permessi per modificare la playlist.
#     Inoltre, la sua homepage mostra tale playlist anche se definita come
#     privata." -> "L’utente invitato è ora autorizzato a modificare la
#     playlist che sarà visibile anche nella sua homepage."
$null = $d.Content.Find.Execute(
  "L’utente selezionato ha ora i permessi per modificare la playlist. Inoltre, la sua homepage mostra tale playlist anche se definita come privata.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "L’utente invitato è ora autorizzato a modificare la playlist che sarà visibile anche nella sua homepage.",
  2)
